$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 752.0833
$ws.Range("J19").Value = 814.3333
$ws.Range("L19").Value = 814.3333
$ws.Range("N19").Value = -1164.3333
$ws.Range("H28").Value = 747.5625
$ws.Range("I28").Value = 681.7273
$ws.Range("K28").Value = 681.7273
$ws.Range("M28").Value = -196.7273
$ws.Range("H32").Value = 4000
$ws.Range("I32").Value = 4000
$ws.Range("K32").Value = 4000
$ws.Range("M32").Value = -3674
$ws.Range("H33").Value = 769
$ws.Range("I33").Value = 193.16667
$ws.Range("K33").Value = 193.16667
$ws.Range("M33").Value = 35.83332999999999
$ws.Range("H43").Value = 5219.25
$ws.Range("I43").Value = 4996.5
$ws.Range("K43").Value = 4996.5
$ws.Range("M43").Value = -4927.5
$ws.Range("H100").Value = 1999.5
$ws.Range("I100").Value = 1999.5
$ws.Range("K100").Value = 1999.5
$ws.Range("M100").Value = -1458.5
$ws.Range("H117").Value = 55000
$ws.Range("J117").Value = 55000
$ws.Range("L117").Value = 55000
$ws.Range("N117").Value = -64178
$ws.Range("H132").Value = 1304.0416
$ws.Range("I132").Value = 1339.2174
$ws.Range("J132").Value = 495
$ws.Range("K132").Value = 4017.6522
$ws.Range("L132").Value = 1485
$ws.Range("M132").Value = -1487.6522
$ws.Range("N132").Value = -6545
$ws.Range("H141").Value = 1499.5
$ws.Range("I141").Value = 1499.5
$ws.Range("K141").Value = 4498.5
$ws.Range("M141").Value = 681.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2438.1
$ws.Range("I2").Value = 1147.7142
$ws.Range("K2").Value = 1147.7142
$ws.Range("M2").Value = -1034.7142
$ws.Range("H45").Value = 3649.111
$ws.Range("I45").Value = 1421
$ws.Range("J45").Value = 4285.7144
$ws.Range("K45").Value = 1421
$ws.Range("L45").Value = 4285.7144
$ws.Range("M45").Value = -1044
$ws.Range("N45").Value = -5039.7144
$ws.Range("H61").Value = 11715.857
$ws.Range("I61").Value = 5499.75
$ws.Range("K61").Value = 5499.75
$ws.Range("M61").Value = -5287.75
$ws.Range("H74").Value = 1635
$ws.Range("I74").Value = 1608.9286
$ws.Range("K74").Value = 1608.9286
$ws.Range("M74").Value = -734.9286
$ws.Range("H77").Value = 1635
$ws.Range("I77").Value = 1608.9286
$ws.Range("K77").Value = 8044.643
$ws.Range("M77").Value = -3676.643
$ws.Range("H102").Value = 2255.1667
$ws.Range("I102").Value = 1769.6666
$ws.Range("K102").Value = 1769.6666
$ws.Range("M102").Value = -147.6666
$ws.Range("H116").Value = 2438.1
$ws.Range("I116").Value = 1147.7142
$ws.Range("K116").Value = 1147.7142
$ws.Range("M116").Value = 1146.2858
$ws.Range("H132").Value = 2749
$ws.Range("I132").Value = 2749
$ws.Range("K132").Value = 8247
$ws.Range("M132").Value = -5717
$ws.Range("H136").Value = 11715.857
$ws.Range("I136").Value = 5499.75
$ws.Range("K136").Value = 16499.25
$ws.Range("M136").Value = -13949.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2438.1
$ws.Range("I3").Value = 1147.7142
$ws.Range("K3").Value = 1147.7142
$ws.Range("M3").Value = -1033.7142
$ws.Range("H57").Value = 25000
$ws.Range("J57").Value = 25000
$ws.Range("L57").Value = 25000
$ws.Range("N57").Value = -26440
$ws.Range("H86").Value = 11999
$ws.Range("J86").Value = 11999
$ws.Range("L86").Value = 11999
$ws.Range("N86").Value = -14245
$ws.Range("H89").Value = 11999
$ws.Range("J89").Value = 11999
$ws.Range("L89").Value = 59995
$ws.Range("N89").Value = -71227
$ws.Range("H99").Value = 2370.1428
$ws.Range("I99").Value = 1935.8182
$ws.Range("J99").Value = 3962.6667
$ws.Range("K99").Value = 1935.8182
$ws.Range("L99").Value = 3962.6667
$ws.Range("M99").Value = -437.8181999999999
$ws.Range("N99").Value = -6958.6667
$ws.Range("H123").Value = 45312
$ws.Range("J123").Value = 45312
$ws.Range("L123").Value = 45312
$ws.Range("N123").Value = -55112
$ws.Range("H136").Value = 25000
$ws.Range("J136").Value = 25000
$ws.Range("L136").Value = 25000
$ws.Range("N136").Value = -35200
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1505.8572
$ws.Range("I7").Value = 994.1429000000001
$ws.Range("K7").Value = 994.1429000000001
$ws.Range("M7").Value = -881.1429000000001
$ws.Range("H22").Value = 20400.666
$ws.Range("J22").Value = 502
$ws.Range("L22").Value = 502
$ws.Range("N22").Value = -1202
$ws.Range("H58").Value = 2370.75
$ws.Range("I58").Value = 2567.4443
$ws.Range("K58").Value = 2567.4443
$ws.Range("M58").Value = -2364.4443
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H94").Value = 2486.4167
$ws.Range("J94").Value = 2953
$ws.Range("L94").Value = 2953
$ws.Range("N94").Value = -3855
$ws.Range("H99").Value = 2627.7693
$ws.Range("I99").Value = 2627.7693
$ws.Range("K99").Value = 2627.7693
$ws.Range("M99").Value = -1129.7693
$ws.Range("H105").Value = 11995.25
$ws.Range("J105").Value = 6902.2
$ws.Range("L105").Value = 6902.2
$ws.Range("N105").Value = -10396.2
$ws.Range("H109").Value = 16936.334
$ws.Range("J109").Value = 16936.334
$ws.Range("L109").Value = 16936.334
$ws.Range("N109").Value = -19016.334
$ws.Range("H126").Value = 2627.7693
$ws.Range("I126").Value = 2627.7693
$ws.Range("K126").Value = 7883.3079
$ws.Range("M126").Value = -5413.3079
$ws.Range("H136").Value = 2370.75
$ws.Range("I136").Value = 2567.4443
$ws.Range("K136").Value = 7702.3329
$ws.Range("M136").Value = -5152.3329
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 358.1875
$ws.Range("J12").Value = 402.57144
$ws.Range("L12").Value = 1207.71432
$ws.Range("N12").Value = -1553.71432
$ws.Range("H113").Value = 763.4
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 763.4
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2290.2
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6630.2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5082.5
$ws.Range("I113").Value = 1165
$ws.Range("K113").Value = 1165
$ws.Range("M113").Value = 1005
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 6822.5
$ws.Range("J100").Value = 1300
$ws.Range("L100").Value = 1300
$ws.Range("N100").Value = -2382
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 30001
$ws.Range("J64").Value = 30001
$ws.Range("L64").Value = 30001
$ws.Range("N64").Value = -30497
$ws.Range("H67").Value = 30001
$ws.Range("J67").Value = 30001
$ws.Range("L67").Value = 30001
$ws.Range("N67").Value = -31717
$ws.Range("H126").Value = 2418.6667
$ws.Range("I126").Value = 2418.6667
$ws.Range("K126").Value = 7256.000100000001
$ws.Range("M126").Value = -4786.000100000001
